$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# =========================================================
# Section 1: fix rows 26/27, 50/51, 89/90, 101-104 (reordered)
# =========================================================
# Row 26 <- source data from old row 27
$ws.Range("B26").Value2 = 6732711
$ws.Range("C26").Value2 = "Lithuania A Lyga"
$ws.Range("D26").Value2 = "Lithuania A Lyga"
$ws.Range("E26").Value2 = 45109.58333333334
$ws.Range("F26").Value2 = "Banga Gargzdai"
$ws.Range("G26").Value2 = "FK Zalgiris Vilnius"
$ws.Range("H26").Value2 = 1
$ws.Range("I26").Value2 = 4
$ws.Range("J26").Value2 = "A"
$ws.Range("K26").Value2 = 5
$ws.Range("L26").Value2 = 3.6
$ws.Range("M26").Value2 = 1.571
$ws.Range("N26").Value2 = 11
$ws.Range("O26").Value2 = 4.75
$ws.Range("P26").Value2 = 1.25
$ws.Range("Q26").Value2 = 1.5
$ws.Range("R26").Value2 = 1.975
$ws.Range("S26").Value2 = 1.825
$ws.Range("T26").Value2 = 2.5
$ws.Range("U26").Value2 = 1.8
$ws.Range("V26").Value2 = 2
$ws.Range("W26").Value2 = -1
$ws.Range("X26").Value2 = -1
$ws.Range("Y26").Value2 = 0.25
$ws.Range("Z26").Value2 = -1
$ws.Range("AA26").Value2 = 0.825
$ws.Range("AB26").Value2 = 0.8
$ws.Range("AC26").Value2 = -1

# Row 27 <- source data from old row 26
$ws.Range("B27").Value2 = 6732773
$ws.Range("C27").Value2 = "Lithuania A Lyga"
$ws.Range("D27").Value2 = "Lithuania A Lyga"
$ws.Range("E27").Value2 = 45109.58333333334
$ws.Range("F27").Value2 = "Suduva Marijampole"
$ws.Range("G27").Value2 = "Hegelmann Litauen"
$ws.Range("H27").Value2 = 0
$ws.Range("I27").Value2 = 1
$ws.Range("J27").Value2 = "A"
$ws.Range("K27").Value2 = 5
$ws.Range("L27").Value2 = 3.8
$ws.Range("M27").Value2 = 1.533
$ws.Range("N27").Value2 = 5
$ws.Range("O27").Value2 = 4.2
$ws.Range("P27").Value2 = 1.533
$ws.Range("Q27").Value2 = 1
$ws.Range("R27").Value2 = 1.875
$ws.Range("S27").Value2 = 1.925
$ws.Range("T27").Value2 = 2.5
$ws.Range("U27").Value2 = 1.9
$ws.Range("V27").Value2 = 1.9
$ws.Range("W27").Value2 = -1
$ws.Range("X27").Value2 = -1
$ws.Range("Y27").Value2 = 0.5329999999999999
$ws.Range("Z27").Value2 = 0
$ws.Range("AA27").Value2 = -0
$ws.Range("AB27").Value2 = -1
$ws.Range("AC27").Value2 = 0.8999999999999999

# Row 50 <- source data from old row 51
$ws.Range("B50").Value2 = 6732795
$ws.Range("C50").Value2 = "Lithuania A Lyga"
$ws.Range("D50").Value2 = "Lithuania A Lyga"
$ws.Range("E50").Value2 = 45149.54166666666
$ws.Range("F50").Value2 = "Suduva Marijampole"
$ws.Range("G50").Value2 = "Banga Gargzdai"
$ws.Range("H50").Value2 = 1
$ws.Range("I50").Value2 = 0
$ws.Range("J50").Value2 = "H"
$ws.Range("K50").Value2 = 2.15
$ws.Range("L50").Value2 = 3.2
$ws.Range("M50").Value2 = 3
$ws.Range("N50").Value2 = 2.3
$ws.Range("O50").Value2 = 3.2
$ws.Range("P50").Value2 = 2.7
$ws.Range("Q50").Value2 = -0.25
$ws.Range("R50").Value2 = 2.05
$ws.Range("S50").Value2 = 1.75
$ws.Range("T50").Value2 = 2.25
$ws.Range("U50").Value2 = 1.9
$ws.Range("V50").Value2 = 1.9
$ws.Range("W50").Value2 = 1.3
$ws.Range("X50").Value2 = -1
$ws.Range("Y50").Value2 = -1
$ws.Range("Z50").Value2 = 1.05
$ws.Range("AA50").Value2 = -1
$ws.Range("AB50").Value2 = -1
$ws.Range("AC50").Value2 = 0.8999999999999999

# Row 51 <- source data from old row 50
$ws.Range("B51").Value2 = 6732794
$ws.Range("C51").Value2 = "Lithuania A Lyga"
$ws.Range("D51").Value2 = "Lithuania A Lyga"
$ws.Range("E51").Value2 = 45149.54166666666
$ws.Range("F51").Value2 = "FK Siauliai"
$ws.Range("G51").Value2 = "FK Dziugas Telsiai"
$ws.Range("H51").Value2 = 3
$ws.Range("I51").Value2 = 0
$ws.Range("J51").Value2 = "H"
$ws.Range("K51").Value2 = 1.25
$ws.Range("L51").Value2 = 5
$ws.Range("M51").Value2 = 9
$ws.Range("N51").Value2 = 1.25
$ws.Range("O51").Value2 = 5.25
$ws.Range("P51").Value2 = 9
$ws.Range("Q51").Value2 = -1.75
$ws.Range("R51").Value2 = 2
$ws.Range("S51").Value2 = 1.8
$ws.Range("T51").Value2 = 3
$ws.Range("U51").Value2 = 1.975
$ws.Range("V51").Value2 = 1.825
$ws.Range("W51").Value2 = 0.25
$ws.Range("X51").Value2 = -1
$ws.Range("Y51").Value2 = -1
$ws.Range("Z51").Value2 = 1
$ws.Range("AA51").Value2 = -1
$ws.Range("AB51").Value2 = 0
$ws.Range("AC51").Value2 = -0

# Row 89 <- source data from old row 90
$ws.Range("B89").Value2 = 7326568
$ws.Range("C89").Value2 = "Lithuania A Lyga"
$ws.Range("D89").Value2 = "Lithuania A Lyga"
$ws.Range("E89").Value2 = 45220.375
$ws.Range("F89").Value2 = "Hegelmann Litauen"
$ws.Range("G89").Value2 = "Panevezys"
$ws.Range("H89").Value2 = 0
$ws.Range("I89").Value2 = 0
$ws.Range("J89").Value2 = "D"
$ws.Range("K89").Value2 = 2.375
$ws.Range("L89").Value2 = 3.2
$ws.Range("M89").Value2 = 2.625
$ws.Range("N89").Value2 = 2.7
$ws.Range("O89").Value2 = 3.2
$ws.Range("P89").Value2 = 2.3
$ws.Range("Q89").Value2 = 0
$ws.Range("R89").Value2 = 2.05
$ws.Range("S89").Value2 = 1.75
$ws.Range("T89").Value2 = 2.25
$ws.Range("U89").Value2 = 1.875
$ws.Range("V89").Value2 = 1.925
$ws.Range("W89").Value2 = -1
$ws.Range("X89").Value2 = 2.2
$ws.Range("Y89").Value2 = -1
$ws.Range("Z89").Value2 = 0
$ws.Range("AA89").Value2 = -0
$ws.Range("AB89").Value2 = -1
$ws.Range("AC89").Value2 = 0.925

# Row 90 <- source data from old row 89
$ws.Range("B90").Value2 = 6732827
$ws.Range("C90").Value2 = "Lithuania A Lyga"
$ws.Range("D90").Value2 = "Lithuania A Lyga"
$ws.Range("E90").Value2 = 45220.375
$ws.Range("F90").Value2 = "FK Dziugas Telsiai"
$ws.Range("G90").Value2 = "FK Kauno Zalgiris"
$ws.Range("H90").Value2 = 0
$ws.Range("I90").Value2 = 2
$ws.Range("J90").Value2 = "A"
$ws.Range("K90").Value2 = 6
$ws.Range("L90").Value2 = 3.9
$ws.Range("M90").Value2 = 1.444
$ws.Range("N90").Value2 = 4.75
$ws.Range("O90").Value2 = 3.6
$ws.Range("P90").Value2 = 1.65
$ws.Range("Q90").Value2 = 0.75
$ws.Range("R90").Value2 = 1.9
$ws.Range("S90").Value2 = 1.9
$ws.Range("T90").Value2 = 2.5
$ws.Range("U90").Value2 = 1.95
$ws.Range("V90").Value2 = 1.85
$ws.Range("W90").Value2 = -1
$ws.Range("X90").Value2 = -1
$ws.Range("Y90").Value2 = 0.6499999999999999
$ws.Range("Z90").Value2 = -1
$ws.Range("AA90").Value2 = 0.8999999999999999
$ws.Range("AB90").Value2 = -1
$ws.Range("AC90").Value2 = 0.8500000000000001

# Row 101 <- source data from old row 103
$ws.Range("B101").Value2 = 7465686
$ws.Range("C101").Value2 = "Lithuania A Lyga"
$ws.Range("D101").Value2 = "Lithuania A Lyga"
$ws.Range("E101").Value2 = 45242.41319444445
$ws.Range("F101").Value2 = "FK Kauno Zalgiris"
$ws.Range("G101").Value2 = "Hegelmann Litauen"
$ws.Range("H101").Value2 = 4
$ws.Range("I101").Value2 = 2
$ws.Range("J101").Value2 = "H"
$ws.Range("K101").Value2 = 2.3
$ws.Range("L101").Value2 = 4
$ws.Range("M101").Value2 = 2.3
$ws.Range("N101").Value2 = 2.55
$ws.Range("O101").Value2 = 4
$ws.Range("P101").Value2 = 2.2
$ws.Range("Q101").Value2 = 0.25
$ws.Range("R101").Value2 = 1.8
$ws.Range("S101").Value2 = 2
$ws.Range("T101").Value2 = 2.75
$ws.Range("U101").Value2 = 1.85
$ws.Range("V101").Value2 = 1.95
$ws.Range("W101").Value2 = 1.55
$ws.Range("X101").Value2 = -1
$ws.Range("Y101").Value2 = -1
$ws.Range("Z101").Value2 = 0.8
$ws.Range("AA101").Value2 = -1
$ws.Range("AB101").Value2 = 0.8500000000000001
$ws.Range("AC101").Value2 = -1

# Row 102 <- source data from old row 101
$ws.Range("B102").Value2 = 6732836
$ws.Range("C102").Value2 = "Lithuania A Lyga"
$ws.Range("D102").Value2 = "Lithuania A Lyga"
$ws.Range("E102").Value2 = 45242.41319444445
$ws.Range("F102").Value2 = "FK Siauliai"
$ws.Range("G102").Value2 = "Banga Gargzdai"
$ws.Range("H102").Value2 = 3
$ws.Range("I102").Value2 = 0
$ws.Range("J102").Value2 = "H"
$ws.Range("K102").Value2 = 1.222
$ws.Range("L102").Value2 = 5.5
$ws.Range("M102").Value2 = 9
$ws.Range("N102").Value2 = 1.363
$ws.Range("O102").Value2 = 4.5
$ws.Range("P102").Value2 = 7
$ws.Range("Q102").Value2 = -1.25
$ws.Range("R102").Value2 = 1.9
$ws.Range("S102").Value2 = 1.9
$ws.Range("T102").Value2 = 2.5
$ws.Range("U102").Value2 = 1.975
$ws.Range("V102").Value2 = 1.825
$ws.Range("W102").Value2 = 0.363
$ws.Range("X102").Value2 = -1
$ws.Range("Y102").Value2 = -1
$ws.Range("Z102").Value2 = 0.8999999999999999
$ws.Range("AA102").Value2 = -1
$ws.Range("AB102").Value2 = 0.9750000000000001
$ws.Range("AC102").Value2 = -1

# Row 103 <- source data from old row 104
$ws.Range("B103").Value2 = 6732727
$ws.Range("C103").Value2 = "Lithuania A Lyga"
$ws.Range("D103").Value2 = "Lithuania A Lyga"
$ws.Range("E103").Value2 = 45242.41319444445
$ws.Range("F103").Value2 = "FK Zalgiris Vilnius"
$ws.Range("G103").Value2 = "FK Dainava Alytus"
$ws.Range("H103").Value2 = 1
$ws.Range("I103").Value2 = 0
$ws.Range("J103").Value2 = "H"
$ws.Range("K103").Value2 = 1.285
$ws.Range("L103").Value2 = 5.5
$ws.Range("M103").Value2 = 6.5
$ws.Range("N103").Value2 = 1.3
$ws.Range("O103").Value2 = 5.5
$ws.Range("P103").Value2 = 6
$ws.Range("Q103").Value2 = -1.5
$ws.Range("R103").Value2 = 1.9
$ws.Range("S103").Value2 = 1.9
$ws.Range("T103").Value2 = 2.75
$ws.Range("U103").Value2 = 1.8
$ws.Range("V103").Value2 = 2
$ws.Range("W103").Value2 = 0.3
$ws.Range("X103").Value2 = -1
$ws.Range("Y103").Value2 = -1
$ws.Range("Z103").Value2 = -1
$ws.Range("AA103").Value2 = 0.8999999999999999
$ws.Range("AB103").Value2 = -1
$ws.Range("AC103").Value2 = 1

# Row 104 <- source data from old row 102
$ws.Range("B104").Value2 = 6732837
$ws.Range("C104").Value2 = "Lithuania A Lyga"
$ws.Range("D104").Value2 = "Lithuania A Lyga"
$ws.Range("E104").Value2 = 45242.41319444445
$ws.Range("F104").Value2 = "Suduva Marijampole"
$ws.Range("G104").Value2 = "FK Riteriai"
$ws.Range("H104").Value2 = 0
$ws.Range("I104").Value2 = 3
$ws.Range("J104").Value2 = "A"
$ws.Range("K104").Value2 = 3.6
$ws.Range("L104").Value2 = 3.6
$ws.Range("M104").Value2 = 1.8
$ws.Range("N104").Value2 = 3
$ws.Range("O104").Value2 = 3.6
$ws.Range("P104").Value2 = 2
$ws.Range("Q104").Value2 = 0.25
$ws.Range("R104").Value2 = 2
$ws.Range("S104").Value2 = 1.8
$ws.Range("T104").Value2 = 2.5
$ws.Range("U104").Value2 = 1.975
$ws.Range("V104").Value2 = 1.825
$ws.Range("W104").Value2 = -1
$ws.Range("X104").Value2 = -1
$ws.Range("Y104").Value2 = 1
$ws.Range("Z104").Value2 = -1
$ws.Range("AA104").Value2 = 0.8
$ws.Range("AB104").Value2 = 0.9750000000000001
$ws.Range("AC104").Value2 = -1

# =========================================================
# Section 2: append new rows 110, 111, 112 (new upcoming matches)
# =========================================================

# Copy style (bold/border for col A, date format for col E) from the last existing row (109)
$ws.Range("A109").Copy($ws.Range("A110"))
$ws.Range("E109").Copy($ws.Range("E110"))
$ws.Range("A109").Copy($ws.Range("A111"))
$ws.Range("E109").Copy($ws.Range("E111"))
$ws.Range("A109").Copy($ws.Range("A112"))
$ws.Range("E109").Copy($ws.Range("E112"))

# Row 110
$ws.Range("A110").Value2 = 108
$ws.Range("B110").Value2 = 7862907
$ws.Range("C110").Value2 = "Lithuania A Lyga"
$ws.Range("D110").Value2 = "Lithuania A Lyga"
$ws.Range("E110").Value2 = 45361.33333333334
$ws.Range("F110").Value2 = "FK Siauliai"
$ws.Range("G110").Value2 = "Banga Gargzdai"
$ws.Range("K110").Value2 = 1.4
$ws.Range("L110").Value2 = 4.5
$ws.Range("M110").Value2 = 5.5
$ws.Range("N110").Value2 = 1.3
$ws.Range("O110").Value2 = 5
$ws.Range("P110").Value2 = 6
$ws.Range("Q110").Value2 = -1.5
$ws.Range("R110").Value2 = 2
$ws.Range("S110").Value2 = 1.8
$ws.Range("T110").Value2 = 2.5
$ws.Range("U110").Value2 = 1.8
$ws.Range("V110").Value2 = 2
$ws.Range("W110").Value2 = 0
$ws.Range("X110").Value2 = 0
$ws.Range("Y110").Value2 = 0
$ws.Range("Z110").Value2 = 0
$ws.Range("AA110").Value2 = 0

# Row 111
$ws.Range("A111").Value2 = 109
$ws.Range("B111").Value2 = 7862908
$ws.Range("C111").Value2 = "Lithuania A Lyga"
$ws.Range("D111").Value2 = "Lithuania A Lyga"
$ws.Range("E111").Value2 = 45361.41666666666
$ws.Range("F111").Value2 = "Panevezys"
$ws.Range("G111").Value2 = "FK Dainava Alytus"
$ws.Range("K111").Value2 = 1.2
$ws.Range("L111").Value2 = 5.5
$ws.Range("M111").Value2 = 10
$ws.Range("N111").Value2 = 1.363
$ws.Range("O111").Value2 = 4.5
$ws.Range("P111").Value2 = 6
$ws.Range("Q111").Value2 = -1.25
$ws.Range("R111").Value2 = 1.95
$ws.Range("S111").Value2 = 1.85
$ws.Range("T111").Value2 = 2.25
$ws.Range("U111").Value2 = 1.8
$ws.Range("V111").Value2 = 2
$ws.Range("W111").Value2 = 0
$ws.Range("X111").Value2 = 0
$ws.Range("Y111").Value2 = 0
$ws.Range("Z111").Value2 = 0
$ws.Range("AA111").Value2 = 0

# Row 112
$ws.Range("A112").Value2 = 110
$ws.Range("B112").Value2 = 7862035
$ws.Range("C112").Value2 = "Lithuania A Lyga"
$ws.Range("D112").Value2 = "Lithuania A Lyga"
$ws.Range("E112").Value2 = 45361.5625
$ws.Range("F112").Value2 = "FK Zalgiris Vilnius"
$ws.Range("G112").Value2 = "Suduva Marijampole"
$ws.Range("K112").Value2 = 1.3
$ws.Range("L112").Value2 = 5
$ws.Range("M112").Value2 = 7
$ws.Range("N112").Value2 = 1.285
$ws.Range("O112").Value2 = 5
$ws.Range("P112").Value2 = 7.5
$ws.Range("Q112").Value2 = -1.5
$ws.Range("R112").Value2 = 1.925
$ws.Range("S112").Value2 = 1.875
$ws.Range("T112").Value2 = 2.75
$ws.Range("U112").Value2 = 1.975
$ws.Range("V112").Value2 = 1.825
$ws.Range("W112").Value2 = 0
$ws.Range("X112").Value2 = 0
$ws.Range("Y112").Value2 = 0
$ws.Range("Z112").Value2 = 0
$ws.Range("AA112").Value2 = 0

Write-Host "All edits applied."
